$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "Acetone ketones ppm"
$ws.Range("F2").Value = 35
$ws.Columns.Item(6).ColumnWidth = 20.166666666666668
$ws.Range("F2").Select()
